$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = @'
https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109864/1750183612496/A5MTrC4t3W-camera-video-segment-driver-1750183614996.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSICANEIFG%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T020004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEK7%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCwKeVwP%2FGu7v4ZCQ8vbQkTsA19DXNXA1AfRsB1cyFf8gIgRGPjutNZpYth1J%2B0r9SbqcsL%2BKhYB2sCX%2FJx8WbQBwEq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDIDKaAq0mTZeX8o9liq6AwSOtY6k5FVe%2B7loK%2FPc9au3nnNY2lFAPUvCnIeWr35TnIYRmiOzZMM8xZo2nF0uDzzgLrcLOm7unuR4Y6od8gHJvnnFq4wDhpvyaNu3pMWdTta88x%2FDlhKIWQCnSsyZzfSPigIWWuPxrfwvj3NC5m9qLtXSscaKGpBZ765aqKZMlDpppcfskp5zinbDXOkTvTSqa6XEtJ7d8tljPaV8ERQTIBdn02CJtVlvX0OYZz6CAyY2Dyv6pkjiHN6eEKxm0N5%2FARy6k0vPwg5IjzgQL7yvSvpClrSHi9CtjBKsMHhfTC5%2BNgZw14a8NYjyLAJPD7vDXzYc%2BPRDXK2ZAWnCNlliD4HykRq5FHgamNVF9etjNWTLGcl7PTvZpV9z5U%2BgZ3tH2YJtsbPVDfhL7KkuIhg3GWcIimtesSSCgbFsF7Z3WcKjKnOuTvgR7Jnve%2Bzau9jIdb2mHy5TXg%2FDeFn8SnWUPd1SYgm78pfepln8ymy04qCKBrqyj6YfPdix8AA5qTeQZgB1VGmV6TnmfIHlDQx85JHOMPX3uAIe%2FHeBCP2U34iq%2BvHGfJJ1aAyCimvrC0ACUHB1uxZuq1kw3NrMwgY6pQE0Zu9E1t%2BEvap9tP34BFqVD3%2Bf78Gcc1vd%2FPmUTOw6Dk4r1wiRZiKpPUPKHuYu5IyAaEokgVeBHR08lVWFez2pqzjf2gaVKMwEenIZgDUsWpBsHMKqGBkoDdfcuqQ6n0eWDyCB1w%2BUjNEZRGpNIG4QHQdaiXVmCfsz5a2Tg0gHM0C3xtSfYL3XgKomrtsKYcCoknLb7E6R0Gk%2ByDP4aunKka%2BocdI%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2010%3A00%3A04%20GMT&X-Amz-Signature=04260382d3678fbced78859ff64bb627b68e0da3630f64c5a70d5c5d9529c182
'@

$ws.Range("L2").Value = @'
https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1750183612496/TAdN8NokD3-camera-video-segment-1750183614996.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSICANEIFG%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T020004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEK7%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCwKeVwP%2FGu7v4ZCQ8vbQkTsA19DXNXA1AfRsB1cyFf8gIgRGPjutNZpYth1J%2B0r9SbqcsL%2BKhYB2sCX%2FJx8WbQBwEq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDIDKaAq0mTZeX8o9liq6AwSOtY6k5FVe%2B7loK%2FPc9au3nnNY2lFAPUvCnIeWr35TnIYRmiOzZMM8xZo2nF0uDzzgLrcLOm7unuR4Y6od8gHJvnnFq4wDhpvyaNu3pMWdTta88x%2FDlhKIWQCnSsyZzfSPigIWWuPxrfwvj3NC5m9qLtXSscaKGpBZ765aqKZMlDpppcfskp5zinbDXOkTvTSqa6XEtJ7d8tljPaV8ERQTIBdn02CJtVlvX0OYZz6CAyY2Dyv6pkjiHN6eEKxm0N5%2FARy6k0vPwg5IjzgQL7yvSvpClrSHi9CtjBKsMHhfTC5%2BNgZw14a8NYjyLAJPD7vDXzYc%2BPRDXK2ZAWnCNlliD4HykRq5FHgamNVF9etjNWTLGcl7PTvZpV9z5U%2BgZ3tH2YJtsbPVDfhL7KkuIhg3GWcIimtesSSCgbFsF7Z3WcKjKnOuTvgR7Jnve%2Bzau9jIdb2mHy5TXg%2FDeFn8SnWUPd1SYgm78pfepln8ymy04qCKBrqyj6YfPdix8AA5qTeQZgB1VGmV6TnmfIHlDQx85JHOMPX3uAIe%2FHeBCP2U34iq%2BvHGfJJ1aAyCimvrC0ACUHB1uxZuq1kw3NrMwgY6pQE0Zu9E1t%2BEvap9tP34BFqVD3%2Bf78Gcc1vd%2FPmUTOw6Dk4r1wiRZiKpPUPKHuYu5IyAaEokgVeBHR08lVWFez2pqzjf2gaVKMwEenIZgDUsWpBsHMKqGBkoDdfcuqQ6n0eWDyCB1w%2BUjNEZRGpNIG4QHQdaiXVmCfsz5a2Tg0gHM0C3xtSfYL3XgKomrtsKYcCoknLb7E6R0Gk%2ByDP4aunKka%2BocdI%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2010%3A00%3A04%20GMT&X-Amz-Signature=707b57e01efdce5ba4d5778d04f804a78bdc9e8bb32ac244d7c0c5be382b5b8c
'@

$ws.Range("K3").Value = @'
https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1750182961683/iBjUYEdqdu-camera-video-segment-driver-1750182964183.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSICANEIFG%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T020004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEK7%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCwKeVwP%2FGu7v4ZCQ8vbQkTsA19DXNXA1AfRsB1cyFf8gIgRGPjutNZpYth1J%2B0r9SbqcsL%2BKhYB2sCX%2FJx8WbQBwEq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDIDKaAq0mTZeX8o9liq6AwSOtY6k5FVe%2B7loK%2FPc9au3nnNY2lFAPUvCnIeWr35TnIYRmiOzZMM8xZo2nF0uDzzgLrcLOm7unuR4Y6od8gHJvnnFq4wDhpvyaNu3pMWdTta88x%2FDlhKIWQCnSsyZzfSPigIWWuPxrfwvj3NC5m9qLtXSscaKGpBZ765aqKZMlDpppcfskp5zinbDXOkTvTSqa6XEtJ7d8tljPaV8ERQTIBdn02CJtVlvX0OYZz6CAyY2Dyv6pkjiHN6eEKxm0N5%2FARy6k0vPwg5IjzgQL7yvSvpClrSHi9CtjBKsMHhfTC5%2BNgZw14a8NYjyLAJPD7vDXzYc%2BPRDXK2ZAWnCNlliD4HykRq5FHgamNVF9etjNWTLGcl7PTvZpV9z5U%2BgZ3tH2YJtsbPVDfhL7KkuIhg3GWcIimtesSSCgbFsF7Z3WcKjKnOuTvgR7Jnve%2Bzau9jIdb2mHy5TXg%2FDeFn8SnWUPd1SYgm78pfepln8ymy04qCKBrqyj6YfPdix8AA5qTeQZgB1VGmV6TnmfIHlDQx85JHOMPX3uAIe%2FHeBCP2U34iq%2BvHGfJJ1aAyCimvrC0ACUHB1uxZuq1kw3NrMwgY6pQE0Zu9E1t%2BEvap9tP34BFqVD3%2Bf78Gcc1vd%2FPmUTOw6Dk4r1wiRZiKpPUPKHuYu5IyAaEokgVeBHR08lVWFez2pqzjf2gaVKMwEenIZgDUsWpBsHMKqGBkoDdfcuqQ6n0eWDyCB1w%2BUjNEZRGpNIG4QHQdaiXVmCfsz5a2Tg0gHM0C3xtSfYL3XgKomrtsKYcCoknLb7E6R0Gk%2ByDP4aunKka%2BocdI%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2010%3A00%3A04%20GMT&X-Amz-Signature=7104df388441597005aaef264bf803b160670fa9d049454e6982532a5d4581c8
'@

$ws.Range("L3").Value = @'
https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1750182961683/Foph7ZSROm-camera-video-segment-1750182964183.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSICANEIFG%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T020004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEK7%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCwKeVwP%2FGu7v4ZCQ8vbQkTsA19DXNXA1AfRsB1cyFf8gIgRGPjutNZpYth1J%2B0r9SbqcsL%2BKhYB2sCX%2FJx8WbQBwEq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDIDKaAq0mTZeX8o9liq6AwSOtY6k5FVe%2B7loK%2FPc9au3nnNY2lFAPUvCnIeWr35TnIYRmiOzZMM8xZo2nF0uDzzgLrcLOm7unuR4Y6od8gHJvnnFq4wDhpvyaNu3pMWdTta88x%2FDlhKIWQCnSsyZzfSPigIWWuPxrfwvj3NC5m9qLtXSscaKGpBZ765aqKZMlDpppcfskp5zinbDXOkTvTSqa6XEtJ7d8tljPaV8ERQTIBdn02CJtVlvX0OYZz6CAyY2Dyv6pkjiHN6eEKxm0N5%2FARy6k0vPwg5IjzgQL7yvSvpClrSHi9CtjBKsMHhfTC5%2BNgZw14a8NYjyLAJPD7vDXzYc%2BPRDXK2ZAWnCNlliD4HykRq5FHgamNVF9etjNWTLGcl7PTvZpV9z5U%2BgZ3tH2YJtsbPVDfhL7KkuIhg3GWcIimtesSSCgbFsF7Z3WcKjKnOuTvgR7Jnve%2Bzau9jIdb2mHy5TXg%2FDeFn8SnWUPd1SYgm78pfepln8ymy04qCKBrqyj6YfPdix8AA5qTeQZgB1VGmV6TnmfIHlDQx85JHOMPX3uAIe%2FHeBCP2U34iq%2BvHGfJJ1aAyCimvrC0ACUHB1uxZuq1kw3NrMwgY6pQE0Zu9E1t%2BEvap9tP34BFqVD3%2Bf78Gcc1vd%2FPmUTOw6Dk4r1wiRZiKpPUPKHuYu5IyAaEokgVeBHR08lVWFez2pqzjf2gaVKMwEenIZgDUsWpBsHMKqGBkoDdfcuqQ6n0eWDyCB1w%2BUjNEZRGpNIG4QHQdaiXVmCfsz5a2Tg0gHM0C3xtSfYL3XgKomrtsKYcCoknLb7E6R0Gk%2ByDP4aunKka%2BocdI%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2010%3A00%3A04%20GMT&X-Amz-Signature=f9d7060a2f53075ebe29f0a26d95638be6484b0e835d7bc8b7b1326fa10c2f8c
'@
